$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the existing "date" number-format (already used by A5) to the
# new date cells A6:A10 by copying formatting only, so we reuse the
# existing style entry instead of minting a new numFmt. ---
$ws.Range("A5").Copy()
$ws.Range("A6:A10").PasteSpecial(-4122)

# --- Row 5: fill in the two new columns (How Long / Approval) ---
$ws.Range("B6").Value = "Online & Initialize Disks 1-4"
$ws.Range("D5").Value = "N/A"
$ws.Range("G5").Value = "ES"
$ws.Range("D2").Value = "Local"

# --- New data rows 6-10 ---
$ws.Range("A6").Value = 42754
$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "Admin"
$ws.Range("F6").Value = "Evan"
$ws.Range("G6").Value = "ES"

$ws.Range("A7").Value = 42754
$ws.Range("B7").Value = "Create Spanned Volume on Disks 1-4"
$ws.Range("C7").Value = "No"
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = "Admin"
$ws.Range("F7").Value = "Evan"
$ws.Range("G7").Value = "ES"

$ws.Range("A8").Value = 42754
$ws.Range("B8").Value = "Create Mirrored Volume on Disks 1&2"
$ws.Range("C8").Value = "No"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "Admin"
$ws.Range("F8").Value = "Evan"
$ws.Range("G8").Value = "ES"

$ws.Range("A9").Value = 42754
$ws.Range("B9").Value = "Create RAID-5 Volume Disks 1-4"
$ws.Range("C9").Value = "No"
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "Admin"
$ws.Range("F9").Value = "Evan"
$ws.Range("G9").Value = "ES"

$ws.Range("A10").Value = 42754
$ws.Range("B10").Value = "Create Mounted Volume @ C:/Mount"
$ws.Range("C10").Value = "No"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = "Admin"
$ws.Range("F10").Value = "Evan"
$ws.Range("G10").Value = "ES"

# --- Move the active selection to A11, matching the saved view state ---
$ws.Range("A11").Select()
